$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive updated loading_percent values in this edit
# (A, E, L, N are untouched index/zero columns)
$cols = @("B", "C", "D", "F", "G", "H", "I", "J", "K", "M", "O")

# Each entry: row number followed by the 11 new values in column order
$data = @(
    @(2, 8.206426242576278, 4.863483411788468, 11.21403608622639, 30.73165030247913, 29.35521344251259, 14.58005122708005, 21.13998314545314, 11.49171140042311, 8.067221664488711, 15.47126241210389, 22.245556022343),
    @(3, 7.876910275122402, 4.695234419308715, 11.16588721719967, 30.81586903234069, 29.48024257965406, 14.63044173731923, 21.23930357929753, 11.51304148670034, 7.843058258827648, 15.3736111673862, 22.33488485598711),
    @(4, 7.667855336319391, 4.588127871527969, 11.13820383055632, 30.87450130273371, 29.56553701844482, 14.6634650932043, 21.30402097458916, 11.52783204526006, 7.702619107449617, 15.31554616830352, 22.3940108105435),
    @(5, 7.581089853974728, 4.543570433997799, 11.12740412897781, 30.90013177353175, 29.6024312303738, 14.67744671915179, 21.33133359212523, 11.53428529420209, 7.644763401607969, 15.29237900984571, 22.41918003728186),
    @(6, 7.566591201137371, 4.536118119501563, 11.12564017155981, 30.90449252747695, 29.60868626174517, 14.67980004000876, 21.33592561285201, 11.53538258007726, 7.635121129224427, 15.28856253848498, 22.4234242693274),
    @(7, 7.666691392050192, 4.587530578305341, 11.13805622133466, 30.87483993415671, 29.56602594858967, 14.66365153065395, 21.30438551605992, 11.52791735122305, 7.701841270109965, 15.31523170027716, 22.39434590019886),
    @(8, 8.094272907352794, 4.806285641476564, 11.19704951190493, 30.75925051863781, 29.39654922701454, 14.59699374475055, 21.17345423449587, 11.49871451098098, 7.990557401918223, 15.43721076135525, 22.27546832678455),
    @(9, 8.874841875238557, 5.203362423760657, 11.32724630834702, 30.58763104553822, 29.13219931228478, 14.48278796847282, 20.94629381516631, 11.4548840341347, 8.53129881155002, 15.6905804739003, 22.07632412770488),
    @(10, 9.407873371701063, 5.473627920081162, 11.43116833391236, 30.49526651373819, 28.97988498070056, 14.4089170912316, 20.79739340874996, 11.43086702879929, 8.966193252156424, 15.88418801328143, 21.9507690879026),
    @(11, 9.640721447796162, 5.591572058001407, 11.48010668338806, 30.46059910077146, 28.91978314590344, 14.37748468030904, 20.73355263820384, 11.42171633649748, 9.163090319819805, 15.97362301517134, 21.89816766136817),
    @(12, 9.727451738839148, 5.635491234234723, 11.49886568920356, 30.44853010023483, 28.89835216379927, 14.36589395963485, 20.70993732323009, 11.41850615230316, 9.236360804344249, 16.00766302929757, 21.87889891700329),
    @(13, 9.708837904622644, 5.626065887634311, 11.49481569746598, 30.45108225855919, 28.90290852251083, 14.36837635406129, 20.71499841033608, 11.41918618654387, 9.220638613260174, 16.00032457980759, 21.88301984656958),
    @(14, 9.647886054674309, 5.595200349413039, 11.48164551706981, 30.45958494729297, 28.91799334532145, 14.37652485218113, 20.73159857416614, 11.42144712391762, 9.169144378455286, 15.97642017412024, 21.89656937117095),
    @(15, 9.610361592835307, 5.576196788999525, 11.47360761245456, 30.46493102688477, 28.92740642426544, 14.38155667352342, 20.74183954917062, 11.42286521227708, 9.13743363865407, 15.96179987057059, 21.90495356499432),
    @(16, 9.392457135405397, 5.465817185397431, 11.42800254761865, 30.49768017772342, 28.98399826467641, 14.41101495377945, 20.80164390238471, 11.4315007428526, 8.953146934118095, 15.87836867690804, 21.9542976491853),
    @(17, 9.256268122258108, 5.396803468376405, 11.40044292333986, 30.51965463388863, 29.02107373452129, 14.42964272508894, 20.83932927802568, 11.43725277759638, 8.837833205307817, 15.82751893487612, 21.9857256484154),
    @(18, 9.177032276018092, 5.356639136485838, 11.38474870188176, 30.53298545426865, 29.04326285171281, 14.44056136292042, 20.86137151656176, 11.44072826173091, 8.770689544718943, 15.7984010919197, 22.00422703091716),
    @(19, 9.150051151428537, 5.34296035491808, 11.37946229572791, 30.53761777557356, 29.05092392103634, 14.44429334621243, 20.86889761373595, 11.4419337006544, 8.747816538380629, 15.78856524154754, 22.01056421121928),
    @(20, 9.270859644571972, 5.404198872768598, 11.40336049438503, 30.51724381264065, 29.01703748330032, 14.42763861050664, 20.83527966503748, 11.43662317444732, 8.850193520689723, 15.83291873067167, 21.98233610783638),
    @(21, 9.665828716470246, 5.604286670108136, 11.48550785284846, 30.45705875260513, 28.91352646230696, 14.37412297374019, 20.72670751274037, 11.42077611365201, 9.184304769339663, 15.98343696828737, 21.89257188999861),
    @(22, 9.915524165877301, 5.730711515907385, 11.54051383563559, 30.42389643819835, 28.85362133286566, 14.34096617610435, 20.65901218369729, 11.41190524594058, 9.395131967285476, 16.08280648429792, 21.83769661119166),
    @(23, 9.783046656601782, 5.663640981752681, 11.51103961101584, 30.44103050371782, 28.88488283990374, 14.35849624641468, 20.69484397929937, 11.41650390430702, 9.283309867121984, 16.0296874886984, 21.86663734220573),
    @(24, 9.264265737421848, 5.400856926788656, 11.40204099183047, 30.51833157269347, 29.01885955078765, 14.4285440187812, 20.83710932254159, 11.43690729295882, 8.844608061398974, 15.83047711883076, 21.98386717042772),
    @(25, 8.670455120499511, 5.099582617556975, 11.29052899210865, 30.62814690599376, 29.19638587668847, 14.51191918284475, 21.00458343996349, 11.46530331972423, 8.38816011005347, 15.62064020696223, 22.12655615750435)
)

foreach ($entry in $data) {
    $rowNum = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $newValue = $entry[$i + 1]
        $ws.Range("$colLetter$rowNum").Value = $newValue
    }
}
